$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to be written with a pandas-style index column (A) and an
# extra blank spacer row (row 4) plus a two-row merged header. The new
# "custom writer" drops the index column/spacer row and writes a flat,
# non-merged 3-row header directly above the data.

# 1) Drop the blank spacer row (old row 4) so the data row moves up.
$ws.Rows(4).Delete()

# 2) Drop the pandas index column (old column A) so data shifts left.
$ws.Columns(1).Delete()

# 3) The old header used merged cells; the new layout repeats the label in
# every spanned cell instead, so unmerge everything first.
$ws.Cells.UnMerge()

# 3b) Drop the old centered/boxed header formatting (incl. its implicit
# vertical="top") so the new header style below starts from a clean slate.
$ws.Cells.ClearFormats()

# 4) Re-write the (now unmerged) header cells so the previously-merged
# labels are repeated across their old span.
$ws.Range("A1").Value = "kunde"
$ws.Range("B1").Value = "kunde"
$ws.Range("C1").Value = "kunde"
$ws.Range("D1").Value = "bestellung"
$ws.Range("E1").Value = "bestellung"

$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "adresse"
$ws.Range("C2").Value = "adresse"
$ws.Range("D2").Value = "id"
$ws.Range("E2").Value = "datum"

$ws.Range("A3").Value = ""
$ws.Range("B3").Value = "strasse"
$ws.Range("C3").Value = "stadt"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""

# 5) Header styling: bold, left-aligned, no border (previously centered
# with a thin box border).
$headerRange = $ws.Range("A1:E3")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4131
$headerRange.Borders.LineStyle = -4142

# 6) Data row stays plain/normal formatting, no border.
$dataRange = $ws.Range("A4:E4")
$dataRange.Font.Bold = $false
$dataRange.Borders.LineStyle = -4142

# 7) Column widths for the new layout. The engine snaps ColumnWidth to a
# 1/6-character pixel grid (Excel-COM "width_chars x 6pt, 5px padding"),
# so these inputs are the closest achievable approximation of the target
# widths (8.7109375 / 11.7109375 / 12.7109375 chars).
$ws.Columns("A").ColumnWidth = 7.75
$ws.Columns("B").ColumnWidth = 10.75
$ws.Columns("C:E").ColumnWidth = 11.75

# 8) Freeze panes below the 3-row header (row 4 is the first scrollable row).
[void]$ws.Range("A4").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
